$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell J1
$ws.Range("J1").Value = "Veterans Use"

# Data rows: "Yes" blocks
$ws.Range("J2:J16").Value = "Yes"
$ws.Range("J39:J40").Value = "Yes"
$ws.Range("J62:J70").Value = "Yes"

# Data rows: "No" blocks
$ws.Range("J17:J38").Value = "No"
$ws.Range("J41:J61").Value = "No"

# Column width
$ws.Range("J1").ColumnWidth = 19.28515625

# Selection
$ws.Range("C1").Select()
